# Automatische test-sync: 2025-08-26 20:07:50
#
# Appends a new "Retour status" log row (row 5) to the "Logs" sheet,
# extends the column conditional-formatting ranges to cover it, and
# bumps the matching category counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

# --- Append the new log entry on row 5 -------------------------------
$logs.Range("A5").Value = "Retour status"
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("D5").Value = "Klantenservice / Opvolging"
$logs.Range("F5").Value = "2025-08-26 20:06:54"
$logs.Range("G5").Value = "Nee"
$logs.Range("H5").Value = "Ja"
$logs.Range("I5").Value = "Nee"
$logs.Range("J5").Value = "Nee"

# --- Extend the per-column conditional formatting to include row 5 ---
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range($col + "2:" + $col + "4")
    $newRange = $logs.Range($col + "2:" + $col + "5")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Update the Dashboard rollup count for the same category ---------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 4
